$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Small text fixes to the German (de) column of the ART dictionary sheet.

# TESTNAME: "Autor:innenquiz" -> "Autor:inquiz"
$ws.Range("B2").Value = "Autor:inquiz"

# PROMPT_PAIRS: "Literat:in" -> "Autor:in"
$ws.Range("B6").Value = 'Welcher der Personen ist eine Autor:in?<br/> Klicken Sie auf den Namen, sie haben {{time_out}} Sekunden Zeit zu antworten.'

# WELCOME: "Test: Autor:innenquiz" -> "Quiz: Autoren und Autorinnen"
$ws.Range("B13").Value = "Quiz: Autoren und Autorinnen"

# FINISHED: "Autor:innenquiz" -> "Autor:inquiz"
$ws.Range("B16").Value = '<h4>Das Autor:inquiz ist nun beendet.</h4> Bitte klicken Sie auf "Weiter", um den nächsten Test zu beginnen.'

# Update the active cell selection to match the final state of the file
$ws.Range("B16").Select()
